$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-format D cells keep their string representation (avoid numeric auto-conversion)
$textCells = @("D2","D3","D5","D6","D8","D9","D11","D13","D15","D16","D18","D20","D21","D22","D24","D28","D32","D35","D36","D38","D39","D40","D42","D43","D45","D46","D47","D48")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.631.86"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.250.60"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "580.47"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "184.18"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "3.248.40"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").Value = "6.57"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "3.799.59"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "27.68"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "67.619.11"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "3.213.06"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").Value = "13.54"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "395.13"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("D22").Value = "7.58"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "71.46"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "9.57"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E31").Value = "  -4.88%  "
$ws.Range("D32").Value = "22.65"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "161.38"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").Value = "1.90"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "26.66"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "6.47"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "40.66"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "2.616.98"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "24.73"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("D48").Value = "334.79"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  -0.55%  "
